$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout: column A's width definition used to overlap column B's
#     (both were covered by one "min=1 max=2" col entry). Touching column B's
#     width splits that into its own independent entry, leaving column A on
#     its own "min=1 max=1" entry, same width as before.
$ws.Columns.Item(2).ColumnWidth = 59.83

# --- Row 10: Objetivos -- replace the mis-pasted teacher name with the real objectives text
$ws.Range("B10").Value = "Apresentar e analisar os conceitos básicos de monitoramento, suas aplicações práticas e as interfaces com os demais instrumentos de Política Ambiental."
$ws.Range("C10").Value = "Apresentar e analisar os conceitos básicos de monitoramento, suas aplicações práticas e as interfaces com os demais instrumentos de Política Ambiental."

# --- Row 13: drop the stray "Programa resumido:" label (now belongs on row 14);
#     B13/C13 switch from the wrong "Semestral" value to the teacher identification data
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C13").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Rows.Item(13).EntireRow.AutoFit()

# --- Row 14: now carries "Programa resumido:" + the real short-syllabus text
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Monitoramento da qualidade ambiental."
$ws.Range("C14").Value = "Monitoramento da qualidade ambiental."
$ws.Rows.Item(14).RowHeight = 60
# B14 is a brand-new cell; copy the normal body style from an existing B-column cell
# (a column-width/style authoring quirk means a fresh cell would otherwise inherit
# the wrong default style)
$ws.Range("B9").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# --- Row 15: just the "Short syllabus:" label, drop the duplicated date value
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: "Programa:" + the real syllabus text
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Conceitos de qualidade ambiental, poluição, padrões de qualidade e de emissão. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. Índices de qualidade. Monitoramento como parte integrante de sistema de gestão ambiental."
$ws.Range("C16").Value = "Conceitos de qualidade ambiental, poluição, padrões de qualidade e de emissão. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. Índices de qualidade. Monitoramento como parte integrante de sistema de gestão ambiental."
$ws.Rows.Item(16).RowHeight = 120
$ws.Range("B9").Copy()
$ws.Range("B16").PasteSpecial(-4122)

# --- Row 17: "Syllabus:" label only
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: "Avaliação:" label only, drop the duplicated teacher name
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).EntireRow.AutoFit()

# --- Row 19: "Método:" label (value columns already hold the right text)
$ws.Range("A19").Value = "Método:"

# --- Row 20: "Critério:" label
$ws.Range("A20").Value = "Critério:"

# --- Row 21: "Norma de recuperação:" label, height shrinks from 120 to 60
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 (new): "Bibliografia:" + full bibliography text
$ws.Range("A22").Value = "Bibliografia:"
$bib = "Porto, R.L.:. org.. Técnicas quantitativas para o gerenciamento de recursos hídricos. ABRH e Editora da Universidade. 1997.`nJames, A. ed., Mathematical models in water pollution control. John Wiley & Sons. 1989. `nMota, S.. Preservação e Conservação de Recursos Hídricos. ABES. 2a. edição. 1995.`nSewell, G.H. Administração e controle de qualidade ambiental. EPU. 1998.`nMacknight, A. Handbook of techniques for aquatic sediments sampling. McGraw Hill 1999. `nLoeb, A. Biological monitoring of aquatic systems. McGraw-Hill. 1998."
$ws.Range("B22").Value = $bib
$ws.Range("C22").Value = $bib
$ws.Rows.Item(22).RowHeight = 120
$ws.Range("B9").Copy()
$ws.Range("B22").PasteSpecial(-4122)
